$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 92) continuing the monthly series in row 91,
# copying A91's formatting (date number format) onto the new date cell.
$ws.Range("A91").Copy($ws.Range("A92"))
$ws.Range("A92").Value = 45474

$ws.Range("B92").Value = 0.625332950081441
$ws.Range("C92").Value = 0.177390208026896
